# Applies the diff described in the commit "print, escapr sequence and
# comments added" to the active document:
#   1. Splits the "...make a existing..." run so "a" is wrapped in
#      spellcheck/grammar proofErr markers.
#   2. Splits the "> print(...)" run so the "print(...)" call is bold.
#   3. Moves <w:lastRenderedPageBreak/> from the "2. sep =>" run to the
#      "will be converted to string before printed." run.
#   4. Splits the "...default one is '<sp>'." run, wrapping the quoted
#      space in gramStart/gramEnd proofErr markers.
#   5. Appends a new " Default is \n." run after "...print at the end."
#   6. Splits the "...Default is sys.stdout." run, wrapping "sys.stdout"
#      in gramStart/gramEnd proofErr markers.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" ' +
  'pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function New-PkgXml([string]$innerParagraphXml) {
    return $pkgHeader + $innerParagraphXml + $pkgFooter
}

function Find-ParagraphByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $para = $doc.Paragraphs($i)
        if ($para.Range.Text -like "*$needle*") {
            return $para
        }
    }
    throw "paragraph containing '$needle' not found"
}

function Set-ParagraphRunXml($doc, [string]$needle, [string]$newParagraphInnerXml) {
    $para = Find-ParagraphByText $doc $needle
    $start = $para.Range.Start
    $end = $para.Range.End - 1   # exclude the paragraph mark
    $range = $doc.Range($start, $end)
    $xml = New-PkgXml("<w:p>" + $newParagraphInnerXml + "</w:p>")
    $range.InsertXML($xml)
}

$lsq = [char]0x2018   # ‘
$rsq = [char]0x2019   # ’

# -- 1. "make a existing" -> "make " + proofErr-wrapped "a" + " existing..."
$p1 = "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>" +
      "<w:t xml:space=`"preserve`">-- To write a comment or make </w:t></w:r>" +
      "<w:proofErr w:type=`"spellStart`"/><w:proofErr w:type=`"gramStart`"/>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:t>a</w:t></w:r>" +
      "<w:proofErr w:type=`"spellEnd`"/><w:proofErr w:type=`"gramEnd`"/>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>" +
      "<w:t xml:space=`"preserve`"> existing single line code as comment, just add $lsq#$rsq at the start of the line.</w:t></w:r>"
Set-ParagraphRunXml $d "a existing single line code" $p1

# -- 2. "> print(...)" -> "> " (plain) + "print(...)" (bold)
$p2 = "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>" +
      "<w:tab/><w:t xml:space=`"preserve`">&gt; </w:t></w:r>" +
      "<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>" +
      "<w:t>print(object(s), sep=separator, end=end, file=file, flush=flush)</w:t></w:r>"
Set-ParagraphRunXml $d "print(object(s), sep=separator" $p2

# -- 3. add <w:lastRenderedPageBreak/> to "will be converted to string before printed."
$p3 = "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>" +
      "<w:lastRenderedPageBreak/>" +
      "<w:t xml:space=`"preserve`">       will be converted to string before printed.</w:t></w:r>"
Set-ParagraphRunXml $d "will be converted to string before printed" $p3

# -- 4. "2. sep =>" (drop lastRenderedPageBreak) + tab + split "default one is '<sp>'."
$p4 = "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:t>2. sep =&gt;</w:t></w:r>" +
      "<w:r w:rsidR=`"00625A4D`"><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:tab/></w:r>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>" +
      "<w:t xml:space=`"preserve`"> Specify how to separate the object, if there is    more than one. The default one is </w:t></w:r>" +
      "<w:proofErr w:type=`"gramStart`"/>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:t>$lsq $lsq</w:t></w:r>" +
      "<w:proofErr w:type=`"gramEnd`"/>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:t>.</w:t></w:r>"
Set-ParagraphRunXml $d "Specify how to separate the object" $p4

# -- 5. append " Default is \n." after "...print at the end."
$p5 = "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:t>3. end =&gt;</w:t></w:r>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>" +
      "<w:tab/><w:t>Specify what to print at the end.</w:t></w:r>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>" +
      "<w:t xml:space=`"preserve`"> Default is \n.</w:t></w:r>"
Set-ParagraphRunXml $d "Specify what to print at the end" $p5

# -- 6. split "An object with a write method. Default is sys.stdout." + trailing tab run
$p6 = "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:t>4. file =&gt;</w:t></w:r>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr>" +
      "<w:tab/><w:t xml:space=`"preserve`">An object with a write method. Default is </w:t></w:r>" +
      "<w:proofErr w:type=`"gramStart`"/>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:t>sys.stdout</w:t></w:r>" +
      "<w:proofErr w:type=`"gramEnd`"/>" +
      "<w:r><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:t>.</w:t></w:r>" +
      "<w:r w:rsidR=`"00625A4D`"><w:rPr><w:sz w:val=`"32`"/><w:szCs w:val=`"32`"/></w:rPr><w:tab/></w:r>"
Set-ParagraphRunXml $d "An object with a write method" $p6

Write-Output "All hunks applied."
